# Add "hydrogen combined cycle" as a new power plant type row, and rename
# the existing "hydrogen" row label to "hydrogen combustion turbine".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Rename the existing row 24 label from "hydrogen" -> "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# New row 25: "hydrogen combined cycle", all-zero guaranteed-dispatch series
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25:AK25").Value = 0

# Match the vertical-center alignment style used on the other fuel-type labels
$ws.Range("A24:A25").VerticalAlignment = -4108

# Mirror the author's final selection/view state
$ws.Range("B25:AK25").Select()
